$d = $word.ActiveDocument

# --- Locate the "Comment out COPY and ENTRYPOINT..." paragraph ---
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Comment out COPY and ENTRYPOINT line in Dockerfile?*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find the 'Comment out COPY and ENTRYPOINT...' paragraph"
}
$targetIndex = $targetPara.Index

# --- Insert a new paragraph right after it ---
$targetPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)

# Type the new sentence, plus a temporary one-character placeholder ("X")
# that we will wrap the _GoBack bookmark around and then delete. Doing the
# bookmark placement this way (rather than adding it to a collapsed Range,
# which this host does not resolve correctly) keeps the sentence in a
# single run and leaves the bookmark collapsed at the very end of the
# paragraph's text, matching how Word itself tracks the last edit
# position ("_GoBack").
$newPara.Range.Text = "Remember to uncomment them if you use the virtual machine.X"

# Re-fetch the paragraph (Range/Paragraph handles can go stale after a
# text write) and build a one-character Range around the placeholder.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$paraEnd = $newPara.Range.End - 1
$markerRange = $d.Range($paraEnd - 1, $paraEnd)

# Move (re-create) the document's _GoBack bookmark onto the placeholder
# character. Word only ever keeps a single _GoBack bookmark, so adding a
# new one here implicitly removes the old one that was on the
# "Note: It might..." paragraph.
$markerRange.Bookmarks.Add("_GoBack")

# Delete the placeholder character; the bookmark collapses to zero length
# right after "machine." (just before the paragraph mark), and the
# sentence's run is left whole again.
$markerRange.Text = ""
